# Weekly update: add new price observations for "Betarraga" (Vega Monumental
# Concepción) and insert them in chronological order among the existing rows,
# shifting the subsequent rows down (matches the row-insert pattern in the
# upstream diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert rows 240:241 (new dates, 2022-10-12 / serial 44846) ----------
$ws.Rows("240:241").Insert()

$ws.Range("A240").Value2 = 11
$ws.Range("B240").Value2 = "Vega Monumental Concepción"
$ws.Range("C240").Value2 = "Bíobío"
$ws.Range("D240").Value2 = 44846
$ws.Range("E240").Value2 = 8
$ws.Range("F240").Value2 = 100114014
$ws.Range("G240").Value2 = "Betarraga"
$ws.Range("H240").Value2 = "Sin especificar"
$ws.Range("I240").Value2 = "Primera"
$ws.Range("J240").Value2 = 250
$ws.Range("K240").Value2 = 800
$ws.Range("L240").Value2 = 800
$ws.Range("M240").Value2 = 800
$ws.Range("N240").Value2 = "$/paquete 5 unidades"
$ws.Range("O240").Value2 = "Región Metropolitana"
$ws.Range("P240").Value2 = 160
$ws.Range("Q240").Value2 = 5
$ws.Range("R240").Value2 = "Hortaliza"

$ws.Range("A241").Value2 = 11
$ws.Range("B241").Value2 = "Vega Monumental Concepción"
$ws.Range("C241").Value2 = "Bíobío"
$ws.Range("D241").Value2 = 44846
$ws.Range("E241").Value2 = 8
$ws.Range("F241").Value2 = 100114014
$ws.Range("G241").Value2 = "Betarraga"
$ws.Range("H241").Value2 = "Sin especificar"
$ws.Range("I241").Value2 = "Segunda"
$ws.Range("J241").Value2 = 300
$ws.Range("K241").Value2 = 650
$ws.Range("L241").Value2 = 650
$ws.Range("M241").Value2 = 650
$ws.Range("N241").Value2 = "$/paquete 5 unidades"
$ws.Range("O241").Value2 = "Región Metropolitana"
$ws.Range("P241").Value2 = 130
$ws.Range("Q241").Value2 = 5
$ws.Range("R241").Value2 = "Hortaliza"

# --- Insert rows 254:255 (new dates, 2022-10-11 / serial 44845) ----------
$ws.Rows("254:255").Insert()

$ws.Range("A254").Value2 = 11
$ws.Range("B254").Value2 = "Vega Monumental Concepción"
$ws.Range("C254").Value2 = "Bíobío"
$ws.Range("D254").Value2 = 44845
$ws.Range("E254").Value2 = 8
$ws.Range("F254").Value2 = 100114014
$ws.Range("G254").Value2 = "Betarraga"
$ws.Range("H254").Value2 = "Sin especificar"
$ws.Range("I254").Value2 = "Primera"
$ws.Range("J254").Value2 = 700
$ws.Range("K254").Value2 = 800
$ws.Range("L254").Value2 = 850
$ws.Range("M254").Value2 = 821
$ws.Range("N254").Value2 = "$/paquete 5 unidades"
$ws.Range("O254").Value2 = "Región Metropolitana"
$ws.Range("P254").Value2 = 164
$ws.Range("Q254").Value2 = 5
$ws.Range("R254").Value2 = "Hortaliza"

$ws.Range("A255").Value2 = 11
$ws.Range("B255").Value2 = "Vega Monumental Concepción"
$ws.Range("C255").Value2 = "Bíobío"
$ws.Range("D255").Value2 = 44845
$ws.Range("E255").Value2 = 8
$ws.Range("F255").Value2 = 100114014
$ws.Range("G255").Value2 = "Betarraga"
$ws.Range("H255").Value2 = "Sin especificar"
$ws.Range("I255").Value2 = "Segunda"
$ws.Range("J255").Value2 = 500
$ws.Range("K255").Value2 = 700
$ws.Range("L255").Value2 = 700
$ws.Range("M255").Value2 = 700
$ws.Range("N255").Value2 = "$/paquete 5 unidades"
$ws.Range("O255").Value2 = "Región Metropolitana"
$ws.Range("P255").Value2 = 140
$ws.Range("Q255").Value2 = 5
$ws.Range("R255").Value2 = "Hortaliza"
